$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A:E (from ~7.71 chars to ~13.38 chars)
$ws.Range("A1:E1").EntireColumn.ColumnWidth = 12.5

# Update the factor-correlation matrix values (keeping it symmetric)
$ws.Range("C2").Value = 0.55
$ws.Range("D2").Value = 0.107
$ws.Range("E2").Value = 0.077

$ws.Range("B3").Value = 0.55
$ws.Range("D3").Value = 0.003
$ws.Range("E3").Value = 0.023

$ws.Range("B4").Value = 0.107
$ws.Range("C4").Value = 0.003
$ws.Range("E4").Value = 0.686

$ws.Range("B5").Value = 0.077
$ws.Range("C5").Value = 0.023
$ws.Range("D5").Value = 0.686
